# "Fixes for Administration page"
#
# 1. Narrow column 1 (600 -> 567 twips) and widen column 2 (7605 -> 7638
#    twips) of the User-Stories table; column 3 (priority) is unchanged.
# 2. Row 2: clarify the save-button note.
# 3. Row 10: user story now talks about saving the process as a .zip
#    instead of opening the evaluated sub-processes.
# 4. Row 12 ("Berechtigungsebenen"): highlight it in yellow instead of
#    green to flag it for follow-up.
# 5. Row 15: fix the "login to identify" story text and drop the
#    "(OpenID?)" aside.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Column widths -----------------------------------------------
# Widths on the Word object model are expressed in points (1 pt = 20
# twips), so 567/20 = 28.35 and 7638/20 = 381.9.
$t.Columns.Item(1).Width = 28.35
$t.Columns.Item(2).Width = 381.9

# --- 2. Row 2 text tweak ---------------------------------------------
$d.Content.Find.Execute("ohne oder mit speichern Botton?)", $true, $false, `
    $false, $false, $false, $true, 1, $false, "mit Speichern-Button)", 2) | Out-Null

# --- 3. Row 10 text tweak ---------------------------------------------
$d.Content.Find.Execute("ebenfalls ausgewertet öffnen.", $true, $false, `
    $false, $false, $false, $true, 1, $false, "als .zip gespeichert werden.", 2) | Out-Null

# --- 4. Row 12 shading: green -> yellow -------------------------------
for ($c = 1; $c -le 3; $c++) {
    $t.Cell(12, $c).Shading.BackgroundPatternColor = 65535   # wdColorYellow (0x00FFFF BGR)
}

# --- 5. Row 15 text tweak ----------------------------------------------
$d.Content.Find.Execute("zu indentifizeren (OpenID?)", $true, $false, `
    $false, $false, $false, $true, 1, $false, "zu identifizieren", 2) | Out-Null
